$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite headers: Expense Date, Expense Value, Description, Expense Type
$ws.Range("A1").Value = "Expense Date"
$ws.Range("B1").Value = "Expense Value"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Expense Type"

# Rewrite data row to match new column order.
# A2 holds a date-looking string ("2024-02-04") that must stay literal TEXT
# (not get auto-converted to a date serial number). Writing it through
# .Formula first (a quoted string literal, never date-inferred), then
# copy/paste-special-values collapses it back down to a plain shared-string
# cell without picking up any numeric/date formatting or quote-prefix style.
$ws.Range("A2").Formula = '="2024-02-04"'
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)
$ws.Range("B2").Value = 100.0
$ws.Range("C2").Value = "New 01"
$ws.Range("D2").Value = "OTHER"

# Remove the now unused column E (old ID/Expense Type leftover column)
$ws.Range("E1:E2").Delete()
